# The numpy-slides code blocks each used to start with a Quarto/Reveal
# caption directive:
#   #| caption: "▶ Ctrl/Cmd+Enter | ⇥ Ctrl/Cmd+] | ⇤ Ctrl/Cmd+["
# followed by a newline and then the real sample code. The site rebuild
# stripped that directive line (plus its trailing newline) from every
# such code block, leaving the rest of each run's text untouched.
#
# We find each occurrence of the directive via TextRange.Find and delete
# exactly that line + trailing newline (60 chars + "\n" = 61 chars),
# repeating per shape for however many times the directive appears there.

$p = $ppt.ActivePresentation

$needle  = '#| caption:'
$lineLen = 61   # length of the caption line (60 chars) + trailing "`n"

# slide index -> list of (shape index, number of occurrences in that shape)
$targets = @(
    @{ Slide = 6;  Shape = 1; Count = 1 },
    @{ Slide = 6;  Shape = 2; Count = 1 },
    @{ Slide = 8;  Shape = 2; Count = 2 },
    @{ Slide = 9;  Shape = 2; Count = 1 },
    @{ Slide = 10; Shape = 2; Count = 2 },
    @{ Slide = 11; Shape = 2; Count = 3 },
    @{ Slide = 12; Shape = 2; Count = 1 },
    @{ Slide = 13; Shape = 2; Count = 2 },
    @{ Slide = 14; Shape = 2; Count = 2 },
    @{ Slide = 15; Shape = 2; Count = 2 },
    @{ Slide = 16; Shape = 2; Count = 1 },
    @{ Slide = 17; Shape = 2; Count = 1 },
    @{ Slide = 18; Shape = 2; Count = 2 },
    @{ Slide = 19; Shape = 2; Count = 1 }
)

foreach ($target in $targets) {
    $s  = $p.Slides.Item($target.Slide)
    $sh = $s.Shapes.Item($target.Shape)
    $tr = $sh.TextFrame.TextRange

    # Bounded loop (never unbounded `while`) so a shape whose text can't
    # be rewritten (e.g. one containing an embedded equation) can't spin
    # forever re-finding the same match.
    for ($n = 0; $n -lt $target.Count; $n++) {
        $found = $tr.Find($needle)
        if ($found -eq $null) {
            break
        }
        $start = $found.Start
        $sub = $tr.Characters($start, $lineLen)
        $sub.Text = ""
    }
}
